$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$passFill = 13434828   # light green interior (matches existing PASSED rows' fill, style index 1/2)
$failFill = 13421823   # light red/pink interior (matches existing FAILED rows' fill, style index 3/4)
$linkFontColor = 15597568   # blue hyperlink font color (matches existing hyperlink cells, style index 2/4)

# ---- Row 56 (test_desktop_anonymous_booking) ----
$ws.Range("A56").Value = "test_desktop_anonymous_booking"
$ws.Range("B56").Value = ""
$ws.Range("C56").Value = ""
$ws.Range("D56").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E56").Value = "FAILED"
$ws.Range("F56").Value = "2025-08-30 20:56:47"
$ws.Range("G56").Value = "46.64s"
$ws.Range("H56").Value = "chrome"
$ws.Range("I56").Value = "Windows"
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""
$ws.Range("L56").Value = ""
$ws.Range("M56").Value = ""
$ws.Range("A56:O56").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N56"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_20-56-47.png", "", "", "📷 Screenshot")
$ws.Range("N56").Font.Color = $linkFontColor
$ws.Range("N56").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O56"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_20-56-47.log", "", "", "🧾 Log File")
$ws.Range("O56").Font.Color = $linkFontColor
$ws.Range("O56").Font.Underline = 2

# ---- Row 57 (test_desktop_anonymous_booking) ----
$ws.Range("A57").Value = "test_desktop_anonymous_booking"
$ws.Range("B57").Value = ""
$ws.Range("C57").Value = ""
$ws.Range("D57").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E57").Value = "FAILED"
$ws.Range("F57").Value = "2025-08-30 20:59:03"
$ws.Range("G57").Value = "45.53s"
$ws.Range("H57").Value = "chrome"
$ws.Range("I57").Value = "Windows"
$ws.Range("J57").Value = ""
$ws.Range("K57").Value = ""
$ws.Range("L57").Value = ""
$ws.Range("M57").Value = ""
$ws.Range("A57:O57").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N57"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_20-59-03.png", "", "", "📷 Screenshot")
$ws.Range("N57").Font.Color = $linkFontColor
$ws.Range("N57").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O57"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_20-59-03.log", "", "", "🧾 Log File")
$ws.Range("O57").Font.Color = $linkFontColor
$ws.Range("O57").Font.Underline = 2

# ---- Row 58 (test_desktop_anonymous_booking) ----
$ws.Range("A58").Value = "test_desktop_anonymous_booking"
$ws.Range("B58").Value = ""
$ws.Range("C58").Value = ""
$ws.Range("D58").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E58").Value = "FAILED"
$ws.Range("F58").Value = "2025-08-30 20:59:56"
$ws.Range("G58").Value = "46.07s"
$ws.Range("H58").Value = "chrome"
$ws.Range("I58").Value = "Windows"
$ws.Range("J58").Value = ""
$ws.Range("K58").Value = ""
$ws.Range("L58").Value = ""
$ws.Range("M58").Value = ""
$ws.Range("A58:O58").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N58"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_20-59-56.png", "", "", "📷 Screenshot")
$ws.Range("N58").Font.Color = $linkFontColor
$ws.Range("N58").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O58"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_20-59-56.log", "", "", "🧾 Log File")
$ws.Range("O58").Font.Color = $linkFontColor
$ws.Range("O58").Font.Underline = 2

# ---- Row 59 (test_desktop_anonymous_booking) ----
$ws.Range("A59").Value = "test_desktop_anonymous_booking"
$ws.Range("B59").Value = ""
$ws.Range("C59").Value = ""
$ws.Range("D59").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E59").Value = "FAILED"
$ws.Range("F59").Value = "2025-08-30 21:02:23"
$ws.Range("G59").Value = "45.85s"
$ws.Range("H59").Value = "chrome"
$ws.Range("I59").Value = "Windows"
$ws.Range("J59").Value = ""
$ws.Range("K59").Value = ""
$ws.Range("L59").Value = ""
$ws.Range("M59").Value = ""
$ws.Range("A59:O59").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N59"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_21-02-23.png", "", "", "📷 Screenshot")
$ws.Range("N59").Font.Color = $linkFontColor
$ws.Range("N59").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O59"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-02-23.log", "", "", "🧾 Log File")
$ws.Range("O59").Font.Color = $linkFontColor
$ws.Range("O59").Font.Underline = 2

# ---- Row 60 (test_desktop_anonymous_booking) ----
$ws.Range("A60").Value = "test_desktop_anonymous_booking"
$ws.Range("B60").Value = ""
$ws.Range("C60").Value = ""
$ws.Range("D60").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E60").Value = "FAILED"
$ws.Range("F60").Value = "2025-08-30 21:05:38"
$ws.Range("G60").Value = "45.38s"
$ws.Range("H60").Value = "chrome"
$ws.Range("I60").Value = "Windows"
$ws.Range("J60").Value = ""
$ws.Range("K60").Value = ""
$ws.Range("L60").Value = ""
$ws.Range("M60").Value = ""
$ws.Range("A60:O60").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N60"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_21-05-38.png", "", "", "📷 Screenshot")
$ws.Range("N60").Font.Color = $linkFontColor
$ws.Range("N60").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O60"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-05-38.log", "", "", "🧾 Log File")
$ws.Range("O60").Font.Color = $linkFontColor
$ws.Range("O60").Font.Underline = 2

# ---- Row 61 (test_desktop_anonymous_booking) ----
$ws.Range("A61").Value = "test_desktop_anonymous_booking"
$ws.Range("B61").Value = ""
$ws.Range("C61").Value = ""
$ws.Range("D61").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E61").Value = "FAILED"
$ws.Range("F61").Value = "2025-08-30 21:07:08"
$ws.Range("G61").Value = "29.67s"
$ws.Range("H61").Value = "chrome"
$ws.Range("I61").Value = "Windows"
$ws.Range("J61").Value = ""
$ws.Range("K61").Value = ""
$ws.Range("L61").Value = ""
$ws.Range("M61").Value = ""
$ws.Range("A61:O61").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N61"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_21-07-08.png", "", "", "📷 Screenshot")
$ws.Range("N61").Font.Color = $linkFontColor
$ws.Range("N61").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O61"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-07-08.log", "", "", "🧾 Log File")
$ws.Range("O61").Font.Color = $linkFontColor
$ws.Range("O61").Font.Underline = 2

# ---- Row 62 (test_desktop_anonymous_booking) ----
$ws.Range("A62").Value = "test_desktop_anonymous_booking"
$ws.Range("B62").Value = ""
$ws.Range("C62").Value = "256898768"
$ws.Range("D62").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E62").Value = "FAILED"
$ws.Range("F62").Value = "2025-08-30 21:08:25"
$ws.Range("G62").Value = "39.39s"
$ws.Range("H62").Value = "chrome"
$ws.Range("I62").Value = "Windows"
$ws.Range("J62").Value = ""
$ws.Range("K62").Value = ""
$ws.Range("L62").Value = ""
$ws.Range("M62").Value = ""
$ws.Range("A62:O62").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N62"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_21-08-25.png", "", "", "📷 Screenshot")
$ws.Range("N62").Font.Color = $linkFontColor
$ws.Range("N62").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O62"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-08-25.log", "", "", "🧾 Log File")
$ws.Range("O62").Font.Color = $linkFontColor
$ws.Range("O62").Font.Underline = 2

# ---- Row 63 (test_desktop_anonymous_booking) ----
$ws.Range("A63").Value = "test_desktop_anonymous_booking"
$ws.Range("B63").Value = ""
$ws.Range("C63").Value = "888438629"
$ws.Range("D63").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E63").Value = "FAILED"
$ws.Range("F63").Value = "2025-08-30 21:13:41"
$ws.Range("G63").Value = "38.18s"
$ws.Range("H63").Value = "chrome"
$ws.Range("I63").Value = "Windows"
$ws.Range("J63").Value = ""
$ws.Range("K63").Value = ""
$ws.Range("L63").Value = ""
$ws.Range("M63").Value = ""
$ws.Range("A63:O63").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("N63"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_21-13-41.png", "", "", "📷 Screenshot")
$ws.Range("N63").Font.Color = $linkFontColor
$ws.Range("N63").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O63"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-13-41.log", "", "", "🧾 Log File")
$ws.Range("O63").Font.Color = $linkFontColor
$ws.Range("O63").Font.Underline = 2

# ---- Row 64 (test_desktop_anonymous_booking) ----
$ws.Range("A64").Value = "test_desktop_anonymous_booking"
$ws.Range("B64").Value = ""
$ws.Range("C64").Value = "503517476"
$ws.Range("D64").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E64").Value = "FAILED"
$ws.Range("F64").Value = "2025-08-30 21:15:37"
$ws.Range("G64").Value = "52.71s"
$ws.Range("H64").Value = "chrome"
$ws.Range("I64").Value = "Windows"
$ws.Range("J64").Value = "חן טסט"
$ws.Range("K64").Value = "restestfattal@gmail.com"
$ws.Range("N64").Value = ""
$ws.Range("A64:O64").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("L64"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/room_selection_2025-08-30_21-15-35.png", "", "", "📷 Screenshot")
$ws.Range("L64").Font.Color = $linkFontColor
$ws.Range("L64").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("M64"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/payment_stage_2025-08-30_21-15-36.png", "", "", "📷 Screenshot")
$ws.Range("M64").Font.Color = $linkFontColor
$ws.Range("M64").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O64"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-15-37.log", "", "", "🧾 Log File")
$ws.Range("O64").Font.Color = $linkFontColor
$ws.Range("O64").Font.Underline = 2

# ---- Row 65 (test_desktop_booking_anonymous_region_eilat) ----
$ws.Range("A65").Value = "test_desktop_booking_anonymous_region_eilat"
$ws.Range("B65").Value = ""
$ws.Range("C65").Value = "794085886"
$ws.Range("D65").Value = "בדיקת השלמת הזמנה משתמש אנונימי דרך אזור מלונות אילת"
$ws.Range("E65").Value = "FAILED"
$ws.Range("F65").Value = "2025-08-30 21:17:29"
$ws.Range("G65").Value = "106.01s"
$ws.Range("H65").Value = "chrome"
$ws.Range("I65").Value = "Windows"
$ws.Range("J65").Value = "חן טסט"
$ws.Range("K65").Value = "restestfattal@gmail.com"
$ws.Range("A65:O65").Interior.Color = $failFill
$ws.Hyperlinks.Add($ws.Range("L65"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/room_selection_2025-08-30_21-17-27.png", "", "", "📷 Screenshot")
$ws.Range("L65").Font.Color = $linkFontColor
$ws.Range("L65").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("M65"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/payment_stage_2025-08-30_21-17-28.png", "", "", "📷 Screenshot")
$ws.Range("M65").Font.Color = $linkFontColor
$ws.Range("M65").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("N65"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_FAIL_2025-08-30_21-17-29.png", "", "", "📷 Screenshot")
$ws.Range("N65").Font.Color = $linkFontColor
$ws.Range("N65").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O65"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_booking_anonymous_region_eilat_2025-08-30_21-17-29.log", "", "", "🧾 Log File")
$ws.Range("O65").Font.Color = $linkFontColor
$ws.Range("O65").Font.Underline = 2

# ---- Row 66 (test_desktop_anonymous_booking) ----
$ws.Range("A66").Value = "test_desktop_anonymous_booking"
$ws.Range("B66").Value = "980025889"
$ws.Range("C66").Value = "435686373"
$ws.Range("D66").Value = "בדיקת השלמת הזמנה מתשמש אנונימי"
$ws.Range("E66").Value = "PASSED"
$ws.Range("F66").Value = "2025-08-30 21:22:34"
$ws.Range("G66").Value = "69.46s"
$ws.Range("H66").Value = "chrome"
$ws.Range("I66").Value = "Windows"
$ws.Range("J66").Value = "חן טסט"
$ws.Range("K66").Value = "restestfattal@gmail.com"
$ws.Range("A66:O66").Interior.Color = $passFill
$ws.Hyperlinks.Add($ws.Range("L66"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/room_selection_2025-08-30_21-22-32.png", "", "", "📷 Screenshot")
$ws.Range("L66").Font.Color = $linkFontColor
$ws.Range("L66").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("M66"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/payment_stage_2025-08-30_21-22-33.png", "", "", "📷 Screenshot")
$ws.Range("M66").Font.Color = $linkFontColor
$ws.Range("M66").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("N66"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/Screenshots/confirmation_PASS_2025-08-30_21-22-34.png", "", "", "📷 Screenshot")
$ws.Range("N66").Font.Color = $linkFontColor
$ws.Range("N66").Font.Underline = 2
$ws.Hyperlinks.Add($ws.Range("O66"), "file:///C:/Users/Chen Ettedgui/PycharmProjects/FattalAUTO/Fattal_Tests/logs/test_desktop_anonymous_booking_2025-08-30_21-22-34.log", "", "", "🧾 Log File")
$ws.Range("O66").Font.Color = $linkFontColor
$ws.Range("O66").Font.Underline = 2
